$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a brand-new worksheet "2022-Q3" right before "2022-Q1" (position 2)
# ---------------------------------------------------------------------------
$insertBefore = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($insertBefore)
$newSheet.Name = "2022-Q3"

# "2022-Q1" got pushed one slot to the right by the insert above, re-fetch it
# by position so we have a stable handle to use as a formatting template.
$template = $wb.Worksheets.Item(3)

# Clone the look of the existing quarter sheets (bold+border+center header,
# bold index column, plain data cells) onto the new sheet before writing data.
$template.Range("A1:H1").Copy($newSheet.Range("A1:H1"))
$template.Range("A2:H2").Copy($newSheet.Range("A2:H2"))
$template.Range("A2:H2").Copy($newSheet.Range("A3:H3"))
$template.Range("A2:H2").Copy($newSheet.Range("A4:H4"))
$template.Range("A2:H2").Copy($newSheet.Range("A5:H5"))
$template.Range("A2:H2").Copy($newSheet.Range("A6:H6"))

# The fund-code / percentage-looking columns must stay text (leading zeros,
# trailing zeros, etc. would otherwise be silently reinterpreted as numbers).
$newSheet.Range("B2:B6").NumberFormat = "@"
$newSheet.Range("D2:G6").NumberFormat = "@"

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data rows (2022-Q3)
$fundData = @(
    @(0, "161017", "富国中证500指数增强（LOF）", "66.37", "90.18", "0.74", "0.4911", 10),
    @(1, "233009", "大摩多因子精选策略混合",       "6.50",  "83.44", "0.92", "0.0598", 8),
    @(2, "013332", "富国中证500指数增强(LOF)C",    "1.68",  "90.18", "0.74", "0.0124", 10),
    @(3, "009613", "上银中证500指数增强A",         "0.97",  "92.48", "1.05", "0.0102", 5),
    @(4, "009614", "上银中证500指数增强C",         "0.74",  "92.48", "1.05", "0.0078", 5)
)

$r = 2
foreach ($entry in $fundData) {
    $newSheet.Range("A$r").Value = $entry[0]
    $newSheet.Range("B$r").Value = $entry[1]
    $newSheet.Range("C$r").Value = $entry[2]
    $newSheet.Range("D$r").Value = $entry[3]
    $newSheet.Range("E$r").Value = $entry[4]
    $newSheet.Range("F$r").Value = $entry[5]
    $newSheet.Range("G$r").Value = $entry[6]
    $newSheet.Range("H$r").Value = $entry[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert the 2022-Q3 totals at the top
#    and push the older quarters down by one row.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Row 6 doesn't exist yet - clone row 5's formatting onto it first so the new
# index cell (A6) keeps the same bold/centered style as A2:A5.
$summary.Range("A5:D5").Copy($summary.Range("A6:D6"))

$summaryRows = @(
    @("2022-Q3", 5, 0.58),
    @("2022-Q1", 2, 0.1),
    @("2021-Q4", 5, 0.26),
    @("2021-Q3", 4, 0.6899999999999999),
    @("2021-Q1", 4, 0.23)
)

$r = 2
$idx = 0
foreach ($entry in $summaryRows) {
    $summary.Range("A$r").Value = $idx
    $summary.Range("B$r").Value = $entry[0]
    $summary.Range("C$r").Value = $entry[1]
    $summary.Range("D$r").Value = $entry[2]
    $r = $r + 1
    $idx = $idx + 1
}
